$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.420.58"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.05%  "
$ws.Range("D3").Value = "'1.864.69"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.18%  "
$ws.Range("E4").Value = "  -1.64%  "
$ws.Range("D5").Value = "'314.12"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.94%  "
$ws.Range("E6").Value = "  -1.46%  "
$ws.Range("D7").Value = "'0.5062"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.94%  "
$ws.Range("D8").Value = "'0.3890"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.79%  "
$ws.Range("D9").Value = "'0.08292"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.52%  "
$ws.Range("D10").Value = "'42.44"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.97%  "
$ws.Range("E11").Value = "  -0.58%  "
$ws.Range("D12").Value = "'6.168"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.13%  "
$ws.Range("D13").Value = "'1.866.98"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.68%  "
$ws.Range("D14").Value = "'20.25"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.85%  "
$ws.Range("D15").Value = "'7.220"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.15%  "
$ws.Range("E16").Value = "  -1.55%  "
$ws.Range("D17").Value = "'0.00001095"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.99%  "
$ws.Range("D18").Value = "'90.96"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.07%  "
$ws.Range("D19").Value = "'0.06715"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.77%  "
$ws.Range("D20").Value = "'17.54"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.66%  "
$ws.Range("E21").Value = "  -1.41%  "
$ws.Range("D22").Value = "'5.889"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.92%  "
$ws.Range("D23").Value = "'28.462.73"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.08%  "
$ws.Range("D24").Value = "'11.02"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.98%  "
$ws.Range("D25").Value = "'2.195"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Value = "'2.074.16"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.18%  "
$ws.Range("D27").Value = "'157.81"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.36%  "
$ws.Range("D28").Value = "'20.46"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.49%  "
$ws.Range("D29").Value = "'2.402"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.94%  "
$ws.Range("D30").Value = "'125.47"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Value = "'0.1031"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.52%  "
$ws.Range("E32").Value = "  +0.18%  "
$ws.Range("D33").Value = "'5.766"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.08%  "
$ws.Range("D34").Value = "'3.615"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.52%  "
$ws.Range("D35").Value = "'0.02435"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.58%  "
$ws.Range("D36").Value = "'0.06560"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.36%  "
$ws.Range("D37").Value = "'8.925"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.46%  "
$ws.Range("D38").Value = "'0.2150"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.43%  "
$ws.Range("D39").Value = "'5.009"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.00%  "
$ws.Range("D40").Value = "'1.177"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.32%  "
$ws.Range("E41").Value = "  -2.91%  "
$ws.Range("D42").Value = "'0.6326"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.14%  "
$ws.Range("D43").Value = "'11.06"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.14%  "
$ws.Range("E44").Value = "  -1.13%  "
$ws.Range("D45").Value = "'0.5959"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.77%  "
$ws.Range("D46").Value = "'13.04"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.19%  "
$ws.Range("D47").Value = "'3.673"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.32%  "
$ws.Range("D48").Value = "'1.988"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.14%  "
$ws.Range("D49").Value = "'121.87"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.05%  "
$ws.Range("E50").Value = "  -0.03%  "
$ws.Range("E51").Value = "  -5.59%  "
